$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 103
$ws1.Range("G4").Value = 70
$ws1.Range("F5").Value = 155
$ws1.Range("F6").Value = 9324
$ws1.Range("F7").Value = 838
$ws1.Range("F10").Value = 1095
$ws1.Range("F12").Value = 79
$ws1.Range("F15").Value = 395
$ws1.Range("F18").Value = 1234

# Sheet "全部类型" (sheet index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 103
$ws4.Range("G6").Value = 70
$ws4.Range("F7").Value = 155
$ws4.Range("F8").Value = 9324
$ws4.Range("F9").Value = 838
$ws4.Range("F12").Value = 1095
$ws4.Range("F14").Value = 79
$ws4.Range("F17").Value = 395
$ws4.Range("F20").Value = 1234
